$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133, shifting existing rows 133..176 down to 134..177
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with its data
$ws.Cells.Item(133, 1).Value = 7
$ws.Cells.Item(133, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(133, 3).Value = "Ñuble"
$ws.Cells.Item(133, 4).Value = 44524
$ws.Cells.Item(133, 5).Value = 16
$ws.Cells.Item(133, 6).Value = 100112043
$ws.Cells.Item(133, 7).Value = "Pepino ensalada"
$ws.Cells.Item(133, 8).Value = "Sin especificar"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 60
$ws.Cells.Item(133, 11).Value = 8000
$ws.Cells.Item(133, 12).Value = 8500
$ws.Cells.Item(133, 13).Value = 8250
$ws.Cells.Item(133, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(133, 15).Value = "Región del Maule"
$ws.Cells.Item(133, 16).Value = 103
$ws.Cells.Item(133, 17).Value = 80
$ws.Cells.Item(133, 18).Value = "Hortaliza"

# Match the date number-format style used by the other rows' Fecha column (D)
$ws.Cells.Item(133, 4).NumberFormat = $ws.Cells.Item(134, 4).NumberFormat
